$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "ECs" cluster (rows 2-4) is dropped. The old "FAPs" cluster (rows 5-7)
# moves up into rows 2-4, and the old "MuSCs" cluster (rows 8-10) moves up
# into rows 5-7. All numeric (TPM-derived) columns are recomputed with the
# new values from the updated script run.

# --- Row 2: FAPs / Ccl28 / Ccr3 / Inflammatory-Mac ---
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Ccl28"
$ws.Range("C2").Value = "Ccr3"
$ws.Range("D2").Value = "Inflammatory-Mac"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.115543
$ws.Range("H2").Value = 0.346629
$ws.Range("I2").Value = 0.7111155332715143
$ws.Range("J2").Value = 0.7111155332715143
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1808983333333333
$ws.Range("N2").Value = 0.5426949999999999
$ws.Range("O2").Value = 0.09546831801815302
$ws.Range("P2").Value = 0.09546831801815302
$ws.Range("Q2").Value = 0.02090153612833333
$ws.Range("R2").Value = 0.188113825155
$ws.Range("S2").Value = 0.0678890038780134
$ws.Range("T2").Value = 0.0678890038780134

# --- Row 3: FAPs / Ccl28 / Ccr3 / Neutrophils ---
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Ccl28"
$ws.Range("C3").Value = "Ccr3"
$ws.Range("D3").Value = "Neutrophils"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.115543
$ws.Range("H3").Value = 0.346629
$ws.Range("I3").Value = 0.7111155332715143
$ws.Range("J3").Value = 0.7111155332715143
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.572737
$ws.Range("N3").Value = 4.718211
$ws.Range("O3").Value = 0.8300051930177132
$ws.Range("P3").Value = 0.8300051930177132
$ws.Range("Q3").Value = 0.181718751191
$ws.Range("R3").Value = 1.635468760719
$ws.Range("S3").Value = 0.5902295854509173
$ws.Range("T3").Value = 0.5902295854509173

# --- Row 4: FAPs / Ccl28 / Ccr3 / Resolving-Mac ---
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ccl28"
$ws.Range("C4").Value = "Ccr3"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.115543
$ws.Range("H4").Value = 0.346629
$ws.Range("I4").Value = 0.7111155332715143
$ws.Range("J4").Value = 0.7111155332715143
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1412166666666667
$ws.Range("N4").Value = 0.42365
$ws.Range("O4").Value = 0.07452648896413371
$ws.Range("P4").Value = 0.07452648896413369
$ws.Range("Q4").Value = 0.01631659731666667
$ws.Range("R4").Value = 0.14684937585
$ws.Range("S4").Value = 0.05299694394258356
$ws.Range("T4").Value = 0.05299694394258356

# --- Row 5: MuSCs / Ccl28 / Ccr3 / Inflammatory-Mac ---
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Ccl28"
$ws.Range("C5").Value = "Ccr3"
$ws.Range("D5").Value = "Inflammatory-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.04693833333333333
$ws.Range("H5").Value = 0.140815
$ws.Range("I5").Value = 0.2888844667284857
$ws.Range("J5").Value = 0.2888844667284857
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1808983333333333
$ws.Range("N5").Value = 0.5426949999999999
$ws.Range("O5").Value = 0.09546831801815302
$ws.Range("P5").Value = 0.09546831801815302
$ws.Range("Q5").Value = 0.008491066269444443
$ws.Range("R5").Value = 0.07641959642499999
$ws.Range("S5").Value = 0.02757931414013962
$ws.Range("T5").Value = 0.02757931414013962

# --- Row 6: MuSCs / Ccl28 / Ccr3 / Neutrophils ---
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Ccl28"
$ws.Range("C6").Value = "Ccr3"
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.04693833333333333
$ws.Range("H6").Value = 0.140815
$ws.Range("I6").Value = 0.2888844667284857
$ws.Range("J6").Value = 0.2888844667284857
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.572737
$ws.Range("N6").Value = 4.718211
$ws.Range("O6").Value = 0.8300051930177132
$ws.Range("P6").Value = 0.8300051930177132
$ws.Range("Q6").Value = 0.07382165355166667
$ws.Range("R6").Value = 0.664394881965
$ws.Range("S6").Value = 0.2397756075667959
$ws.Range("T6").Value = 0.2397756075667959

# --- Row 7: MuSCs / Ccl28 / Ccr3 / Resolving-Mac ---
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Ccl28"
$ws.Range("C7").Value = "Ccr3"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.04693833333333333
$ws.Range("H7").Value = 0.140815
$ws.Range("I7").Value = 0.2888844667284857
$ws.Range("J7").Value = 0.2888844667284857
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.1412166666666667
$ws.Range("N7").Value = 0.42365
$ws.Range("O7").Value = 0.07452648896413371
$ws.Range("P7").Value = 0.07452648896413369
$ws.Range("Q7").Value = 0.006628474972222223
$ws.Range("R7").Value = 0.05965627475
$ws.Range("S7").Value = 0.02152954502155014
$ws.Range("T7").Value = 0.02152954502155014

# The old "MuSCs" cluster rows (8-10) are no longer present; the data that
# used to occupy them has moved up into rows 5-7 above, so delete the
# now-redundant trailing rows (also shrinks the sheet dimension to A1:T7).
$ws.Range("A8:T10").Delete()
